# Update cryptos price/volume table to reflect the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price values that are NOT number-like (contain multiple dots, subscript digits, etc.) ---
# These can be assigned directly since Excel will keep them as text.
$ws.Range("D2").Value = "67.341.68"
$ws.Range("D3").Value = "2.625.51"
$ws.Range("D9").Value = "2.625.09"
$ws.Range("D15").Value = "3.104.10"
$ws.Range("D17").Value = "67.157.29"
$ws.Range("D18").Value = "2.620.41"
$ws.Range("D46").Value = "0.0₆0296"

# --- Price values that DO look like plain numbers to Excel. ---
# Stage them in a scratch column formatted as Text, then paste-special (values-only)
# into the destination so the destination keeps its original (default) cell style
# while the stored value remains a text string, matching the source data feed.
$ws.Range("ZZ1").NumberFormat = "@"
$ws.Range("ZZ1").Value = "595.06"
$ws.Range("ZZ2").NumberFormat = "@"
$ws.Range("ZZ2").Value = "167.15"
$ws.Range("ZZ3").NumberFormat = "@"
$ws.Range("ZZ3").Value = "27.65"
$ws.Range("ZZ4").NumberFormat = "@"
$ws.Range("ZZ4").Value = "0.0000182"
$ws.Range("ZZ5").NumberFormat = "@"
$ws.Range("ZZ5").Value = "12.09"
$ws.Range("ZZ6").NumberFormat = "@"
$ws.Range("ZZ6").Value = "357.48"
$ws.Range("ZZ7").NumberFormat = "@"
$ws.Range("ZZ7").Value = "4.33"
$ws.Range("ZZ8").NumberFormat = "@"
$ws.Range("ZZ8").Value = "1.93"
$ws.Range("ZZ9").NumberFormat = "@"
$ws.Range("ZZ9").Value = "10.28"
$ws.Range("ZZ10").NumberFormat = "@"
$ws.Range("ZZ10").Value = "69.71"
$ws.Range("ZZ11").NumberFormat = "@"
$ws.Range("ZZ11").Value = "544.66"
$ws.Range("ZZ12").NumberFormat = "@"
$ws.Range("ZZ12").Value = "156.69"
$ws.Range("ZZ13").NumberFormat = "@"
$ws.Range("ZZ13").Value = "5.22"
$ws.Range("ZZ14").NumberFormat = "@"
$ws.Range("ZZ14").Value = "1.81"
$ws.Range("ZZ15").NumberFormat = "@"
$ws.Range("ZZ15").Value = "18.19"
$ws.Range("ZZ16").NumberFormat = "@"
$ws.Range("ZZ16").Value = "152.22"

$ws.Range("ZZ1").Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("ZZ2").Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("ZZ3").Copy()
$ws.Range("D14").PasteSpecial(-4163)
$ws.Range("ZZ4").Copy()
$ws.Range("D16").PasteSpecial(-4163)
$ws.Range("ZZ5").Copy()
$ws.Range("D19").PasteSpecial(-4163)
$ws.Range("ZZ6").Copy()
$ws.Range("D21").PasteSpecial(-4163)
$ws.Range("ZZ7").Copy()
$ws.Range("D22").PasteSpecial(-4163)
$ws.Range("ZZ8").Copy()
$ws.Range("D25").PasteSpecial(-4163)
$ws.Range("ZZ9").Copy()
$ws.Range("D26").PasteSpecial(-4163)
$ws.Range("ZZ10").Copy()
$ws.Range("D27").PasteSpecial(-4163)
$ws.Range("ZZ11").Copy()
$ws.Range("D31").PasteSpecial(-4163)
$ws.Range("ZZ12").Copy()
$ws.Range("D38").PasteSpecial(-4163)
$ws.Range("ZZ13").Copy()
$ws.Range("D41").PasteSpecial(-4163)
$ws.Range("ZZ14").Copy()
$ws.Range("D42").PasteSpecial(-4163)
$ws.Range("ZZ15").Copy()
$ws.Range("D43").PasteSpecial(-4163)
$ws.Range("ZZ16").Copy()
$ws.Range("D47").PasteSpecial(-4163)

$ws.Range("ZZ1:ZZ16").Clear()
$excel.CutCopyMode = 0

# --- Coin name / link / volume(1h) text updates ---
$ws.Range("E2").Value = "  -0.41%  "
$ws.Range("E3").Value = "  -2.09%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("E5").Value = "  -0.91%  "
$ws.Range("E6").Value = "  +0.79%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  -2.54%  "
$ws.Range("E10").Value = "  -2.84%  "
$ws.Range("E11").Value = "  +1.29%  "
$ws.Range("E12").Value = "  +1.38%  "
$ws.Range("E13").Value = "  +0.32%  "
$ws.Range("E14").Value = "  -0.70%  "
$ws.Range("E15").Value = "  -2.08%  "
$ws.Range("E16").Value = "  -1.09%  "
$ws.Range("E17").Value = "  -0.68%  "
$ws.Range("E18").Value = "  -1.51%  "
$ws.Range("E19").Value = "  +2.74%  "
$ws.Range("E20").Value = "  +4.11%  "
$ws.Range("E21").Value = "  -1.93%  "
$ws.Range("E22").Value = "  -1.39%  "
$ws.Range("E23").Value = "  -3.26%  "
$ws.Range("E24").Value = "  +0.00%  "
$ws.Range("E25").Value = "  -5.24%  "
$ws.Range("E26").Value = "  +1.57%  "
$ws.Range("E27").Value = "  -2.06%  "
$ws.Range("E29").Value = "  -0.08%  "
$ws.Range("E30").Value = "  -1.98%  "
$ws.Range("E31").Value = "  -2.54%  "
$ws.Range("E32").Value = "  -0.82%  "
$ws.Range("E33").Value = "  -3.17%  "
$ws.Range("E34").Value = "  -1.86%  "
$ws.Range("E35").Value = "  +4.61%  "
$ws.Range("E36").Value = "  +0.07%  "
$ws.Range("E37").Value = "  -2.99%  "
$ws.Range("E38").Value = "  +1.24%  "
$ws.Range("E39").Value = "  -2.86%  "
$ws.Range("E40").Value = "  -2.27%  "
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("E41").Value = "  -1.83%  "
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("E42").Value = "  -1.28%  "
$ws.Range("E43").Value = "  +1.31%  "
$ws.Range("E44").Value = "  +0.08%  "
$ws.Range("E45").Value = "  -4.52%  "
$ws.Range("E46").Value = "  -0.52%  "
$ws.Range("E47").Value = "  -0.83%  "
$ws.Range("E48").Value = "  -1.89%  "
$ws.Range("E49").Value = "  -1.41%  "
$ws.Range("E50").Value = "  -1.41%  "
$ws.Range("E51").Value = "  -0.97%  "
